$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 383; this pushes the existing rows 383-403
# down to 384-404 (carrying all their data/formatting with them, so the
# former row 403 ends up correctly at row 404 with no further edits).
$ws.Rows.Item(383).Insert()

# Populate the newly inserted row 383 with a new weekly price observation.
# Most columns mirror the row immediately below it (the record that used
# to be row 383), only the date/volume/price columns differ.
$ws.Range("A383").Value = 10
$ws.Range("B383").Value = "Vega Modelo de Temuco"
$ws.Range("C383").Value = "La Araucanía"
$ws.Range("D383").Value = "2023-04-05"
$ws.Range("E383").Value = 9
$ws.Range("F383").Value = "Fruta"
$ws.Range("G383").Value = 100102
$ws.Range("H383").Value = "Cítricos"
$ws.Range("I383").Value = 100102006
$ws.Range("J383").Value = "Pomelo"
$ws.Range("K383").Value = "Start Ruby"
$ws.Range("L383").Value = "Primera"
$ws.Range("M383").Value = 130
$ws.Range("N383").Value = 15000
$ws.Range("O383").Value = 16000
$ws.Range("P383").Value = 15385
$ws.Range("Q383").Value = "$/bandeja 15 kilos granel"
$ws.Range("R383").Value = "Región de O'Higgins"
$ws.Range("S383").Value = 1026
$ws.Range("T383").Value = 15
